# Applies the crypto price/volume refresh captured by the source diff.
#
# Rows 2-45 (existing coins): only the Price (D) and/or Volume(1h) (E) text
# labels are refreshed in place.
#
# Rows 46-51: BabyDogeCoin dropped out of the top list, so every row from
# "Aave" onward shifted up by one and "Mantle" was appended as the new last
# row; those six rows get full Coin/Link/Price/Volume(1h) replacements.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> column letter -> new text value.
$rowUpdates = [ordered]@{
    2 = @{ "D" = '26.237.72'; "E" = '  -1.07%  ' }
    3 = @{ "D" = '1.662.03'; "E" = '  -1.11%  ' }
    4 = @{ "E" = '  +0.43%  ' }
    5 = @{ "D" = '218.83'; "E" = '  +1.01%  ' }
    6 = @{ "D" = '0.5225'; "E" = '  -1.82%  ' }
    7 = @{ "E" = '  +0.43%  ' }
    8 = @{ "D" = '0.2671'; "E" = '  -0.36%  ' }
    9 = @{ "D" = '0.06347'; "E" = '  -0.75%  ' }
    10 = @{ "D" = '21.09'; "E" = '  -2.56%  ' }
    11 = @{ "D" = '0.07715'; "E" = '  -1.13%  ' }
    12 = @{ "D" = '1.664.93'; "E" = '  -0.89%  ' }
    13 = @{ "D" = '4.433'; "E" = '  -1.54%  ' }
    14 = @{ "D" = '1.889.25'; "E" = '  -1.00%  ' }
    15 = @{ "D" = '0.5477'; "E" = '  -1.65%  ' }
    16 = @{ "D" = '0.0₅8214'; "E" = '  -1.49%  ' }
    17 = @{ "D" = '65.05'; "E" = '  -1.03%  ' }
    18 = @{ "D" = '26.248.92'; "E" = '  -1.16%  ' }
    19 = @{ "D" = '1.006'; "E" = '  +0.45%  ' }
    20 = @{ "D" = '4.657'; "E" = '  -2.27%  ' }
    21 = @{ "D" = '195.39'; "E" = '  +0.44%  ' }
    22 = @{ "D" = '10.16'; "E" = '  -1.90%  ' }
    23 = @{ "D" = '6.094'; "E" = '  -3.96%  ' }
    24 = @{ "E" = '  +0.61%  ' }
    25 = @{ "D" = '139.27'; "E" = '  -2.78%  ' }
    26 = @{ "D" = '0.1244'; "E" = '  -3.06%  ' }
    27 = @{ "D" = '7.236'; "E" = '  -2.68%  ' }
    28 = @{ "D" = '16.22'; "E" = '  -0.69%  ' }
    29 = @{ "E" = '  -0.80%  ' }
    30 = @{ "D" = '0.05975'; "E" = '  -3.04%  ' }
    31 = @{ "D" = '1.282'; "E" = '  +0.67%  ' }
    32 = @{ "D" = '3.624'; "E" = '  +0.48%  ' }
    33 = @{ "D" = '3.309'; "E" = '  -4.25%  ' }
    34 = @{ "E" = '  -3.35%  ' }
    35 = @{ "D" = '0.9808'; "E" = '  -2.73%  ' }
    36 = @{ "D" = '2.424'; "E" = '  +0.04%  ' }
    37 = @{ "D" = '2.782'; "E" = '  -0.30%  ' }
    38 = @{ "D" = '0.5911'; "E" = '  +3.23%  ' }
    39 = @{ "D" = '0.01597'; "E" = '  -2.44%  ' }
    40 = @{ "D" = '5.996'; "E" = '  -0.50%  ' }
    41 = @{ "D" = '0.8580'; "E" = '  -0.20%  ' }
    43 = @{ "D" = '1.032.04'; "E" = '  -4.00%  ' }
    44 = @{ "D" = '99.94'; "E" = '  -0.13%  ' }
    45 = @{ "D" = '1.803.04'; "E" = '  -1.33%  ' }
    46 = @{ "B" = 'Aave'; "C" = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'; "D" = '57.43'; "E" = '  +0.74%  ' }
    47 = @{ "B" = 'Frax'; "C" = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'; "D" = '1.006'; "E" = '  +0.44%  ' }
    48 = @{ "B" = 'EnergySwap'; "C" = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; "D" = '8.079'; "E" = '  -1.00%  ' }
    49 = @{ "B" = 'Cronos'; "C" = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'; "D" = '0.05187'; "E" = '  -0.44%  ' }
    50 = @{ "B" = 'RenderToken'; "C" = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; "D" = '1.469'; "E" = '  +0.00%  ' }
    51 = @{ "B" = 'Mantle'; "C" = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'; "D" = '0.4229'; "E" = '  -0.25%  ' }
}

foreach ($row in $rowUpdates.Keys) {
    $cols = $rowUpdates[$row]
    foreach ($col in $cols.Keys) {
        $cell = $ws.Range("$col$row")
        # Force text interpretation so numeric-looking labels (e.g. "218.83")
        # are not silently reinterpreted as actual numbers; Style reset drops
        # the transient "@" text format once the literal text is committed, so
        # no cell ends up with a style it did not have before.
        $cell.NumberFormat = "@"
        $cell.Value = $cols[$col]
        $cell.Style = "Normal"
    }
}
